$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at C and D (shifts old full/tipo/link -> E/F/G)
$ws.Columns("C:D").Insert()

# Update header row
$ws.Cells.Item(1, 3).Value2 = "modelo"
$ws.Cells.Item(1, 4).Value2 = "politica"

# Remove old data rows (2-8) so we can rewrite the full data set cleanly
$ws.Rows("2:8").Delete()

# Write the full data set (rows 2-18)
$r = 2
$ws.Cells.Item($r, 1).Value2 = "Fonte Carregador Automotiva Jfa 200a Slim Bivolt Voltímetro"
$ws.Cells.Item($r, 2).Value2 = 845.87
$ws.Cells.Item($r, 3).Value2 = "FONTE 200A"
$ws.Cells.Item($r, 4).Value2 = "Igual"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://www.mercadolivre.com.br/fonte-carregador-automotiva-jfa-200a-slim-bivolt-voltimetro/p/MLB21348561?pdp_filters=seller_id:64243271#searchVariation=MLB21348561&position=2&search_layout=stack&type=product&tracking_id=f745bc25-3749-4911-b299-8e5c83f83893"

$r = 3
$ws.Cells.Item($r, 1).Value2 = "Fonte Automotiva Jfa Storm 200a Bob Carregador Automático Bivolt Cor BOB 200A JFA"
$ws.Cells.Item($r, 2).Value2 = 731.39
$ws.Cells.Item($r, 3).Value2 = "FONTE 200 BOB"
$ws.Cells.Item($r, 4).Value2 = "Acima"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://www.mercadolivre.com.br/fonte-automotiva-jfa-storm-200a-bob-carregador-automatico-bivolt-cor-bob-200a-jfa/p/MLB24834408?pdp_filters=seller_id:64243271#searchVariation=MLB24834408&position=5&search_layout=stack&type=product&tracking_id=f745bc25-3749-4911-b299-8e5c83f83893"

$r = 4
$ws.Cells.Item($r, 1).Value2 = "Fonte Carregador Jfa Bob Storm 90a Bivolt Automático Cor Preto"
$ws.Cells.Item($r, 2).Value2 = 466.39
$ws.Cells.Item($r, 3).Value2 = "FONTE 90 BOB"
$ws.Cells.Item($r, 4).Value2 = "Acima"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://www.mercadolivre.com.br/fonte-carregador-jfa-bob-storm-90a-bivolt-automatico-cor-preto/p/MLB21562641?pdp_filters=seller_id:64243271#searchVariation=MLB21562641&position=4&search_layout=stack&type=product&tracking_id=f745bc25-3749-4911-b299-8e5c83f83893"

$r = 5
$ws.Cells.Item($r, 1).Value2 = "Fonte Carregador Jfa 60a Bivolt Storm Com Medidor Cca"
$ws.Cells.Item($r, 2).Value2 = 533
$ws.Cells.Item($r, 3).Value2 = "Modelo identificado mas fora do range de preco"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://www.mercadolivre.com.br/fonte-carregador-jfa-60a-bivolt-storm-com-medidor-cca/p/MLB21320712?pdp_filters=seller_id:64243271#searchVariation=MLB21320712&position=1&search_layout=stack&type=product&tracking_id=f745bc25-3749-4911-b299-8e5c83f83893"

$r = 6
$ws.Cells.Item($r, 1).Value2 = "Fonte Som Automotiva Jfa Carregador De Bateria 120a Bob"
$ws.Cells.Item($r, 2).Value2 = 539.74
$ws.Cells.Item($r, 3).Value2 = "FONTE 120 BOB"
$ws.Cells.Item($r, 4).Value2 = "Igual"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-3701888097-fonte-som-automotiva-jfa-carregador-de-bateria-120a-bob-_JM#position%3D8%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

$r = 7
$ws.Cells.Item($r, 1).Value2 = "Fonte Som Automotiva Jfa Carregador De Bateria 90a Bob"
$ws.Cells.Item($r, 2).Value2 = 443.07
$ws.Cells.Item($r, 3).Value2 = "FONTE 90 BOB"
$ws.Cells.Item($r, 4).Value2 = "Igual"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-4706615412-fonte-som-automotiva-jfa-carregador-de-bateria-90a-bob-_JM#position%3D9%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

$r = 8
$ws.Cells.Item($r, 1).Value2 = "Fonte Som Automotiva Jfa Carregador De Bateria 200a Bob"
$ws.Cells.Item($r, 2).Value2 = 694.82
$ws.Cells.Item($r, 3).Value2 = "FONTE 200 BOB"
$ws.Cells.Item($r, 4).Value2 = "Igual"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-3701844755-fonte-som-automotiva-jfa-carregador-de-bateria-200a-bob-_JM#position%3D11%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

$r = 9
$ws.Cells.Item($r, 1).Value2 = "Fonte Som Automotiva Jfa Carregador De Bateria 200a Storm"
$ws.Cells.Item($r, 2).Value2 = 845.87
$ws.Cells.Item($r, 3).Value2 = "FONTE 200A"
$ws.Cells.Item($r, 4).Value2 = "Igual"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-3701720815-fonte-som-automotiva-jfa-carregador-de-bateria-200a-storm-_JM#position%3D12%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

$r = 10
$ws.Cells.Item($r, 1).Value2 = "Fonte Som Automotiva Jfa Carregador De Bateria 60a Storm"
$ws.Cells.Item($r, 2).Value2 = 473.28
$ws.Cells.Item($r, 3).Value2 = "FONTE 60A"
$ws.Cells.Item($r, 4).Value2 = "Igual"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-4707095142-fonte-som-automotiva-jfa-carregador-de-bateria-60a-storm-_JM#position%3D13%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

$r = 11
$ws.Cells.Item($r, 1).Value2 = "Volante Esportivo Astra Meriva Montana Zafira Corsa Joy Jfa"
$ws.Cells.Item($r, 2).Value2 = 727.89
$ws.Cells.Item($r, 3).Value2 = "Sem Modelo"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-1803728714-volante-esportivo-astra-meriva-montana-zafira-corsa-joy-jfa-_JM#position%3D14%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

$r = 12
$ws.Cells.Item($r, 1).Value2 = "Volante Esportivo Corsa Classic Kadet Monza Prisma Celta Jfa"
$ws.Cells.Item($r, 2).Value2 = 615.23
$ws.Cells.Item($r, 3).Value2 = "Sem Modelo"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-1803743331-volante-esportivo-corsa-classic-kadet-monza-prisma-celta-jfa-_JM#position%3D15%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

$r = 13
$ws.Cells.Item($r, 1).Value2 = "Volante Esportivo Jetta Gli Montana 2012 2013 Controle Jfa"
$ws.Cells.Item($r, 2).Value2 = 567.81
$ws.Cells.Item($r, 3).Value2 = "Sem Modelo"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-1989056300-volante-esportivo-jetta-gli-montana-2012-2013-controle-jfa-_JM#position%3D16%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

$r = 14
$ws.Cells.Item($r, 1).Value2 = "Volante Esportivo Jetta Gli Astra Meriva Montana Corsa Jfa"
$ws.Cells.Item($r, 2).Value2 = 567.81
$ws.Cells.Item($r, 3).Value2 = "Sem Modelo"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-1989043880-volante-esportivo-jetta-gli-astra-meriva-montana-corsa-jfa-_JM#position%3D17%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

$r = 15
$ws.Cells.Item($r, 1).Value2 = "Volante Esportivo Jetta Gli Astra 98 99 2000 Controle Jfa"
$ws.Cells.Item($r, 2).Value2 = 567.81
$ws.Cells.Item($r, 3).Value2 = "Sem Modelo"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-1989060232-volante-esportivo-jetta-gli-astra-98-99-2000-controle-jfa-_JM#position%3D18%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

$r = 16
$ws.Cells.Item($r, 1).Value2 = "Volante Esportivo Jetta Gli Corsa Hatch 2008 Controle Jfa"
$ws.Cells.Item($r, 2).Value2 = 567.81
$ws.Cells.Item($r, 3).Value2 = "Sem Modelo"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "premium"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-1989051284-volante-esportivo-jetta-gli-corsa-hatch-2008-controle-jfa-_JM#position%3D19%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

$r = 17
$ws.Cells.Item($r, 1).Value2 = "Volante Esportivo Corsa Classic Kadet Monza Prisma Celta Jfa"
$ws.Cells.Item($r, 2).Value2 = 581.23
$ws.Cells.Item($r, 3).Value2 = "Sem Modelo"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "classico"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-1803729492-volante-esportivo-corsa-classic-kadet-monza-prisma-celta-jfa-_JM#position%3D20%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

$r = 18
$ws.Cells.Item($r, 1).Value2 = "Volante Esportivo Astra Meriva Montana Zafira Corsa Joy Jfa"
$ws.Cells.Item($r, 2).Value2 = 690.48
$ws.Cells.Item($r, 3).Value2 = "Sem Modelo"
$ws.Cells.Item($r, 5).Value2 = "NA"
$ws.Cells.Item($r, 6).Value2 = "classico"
$ws.Cells.Item($r, 7).Value2 = "https://produto.mercadolivre.com.br/MLB-1803734688-volante-esportivo-astra-meriva-montana-zafira-corsa-joy-jfa-_JM#position%3D21%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df745bc25-3749-4911-b299-8e5c83f83893"

Write-Output "done"